# Append one new upload-log row (row 46) to each of the four lifter sheets,
# mirroring the existing row layout (time / length / ID / actual length /
# checksum + their decimal counterparts).

function Add-UploadRow {
    param($ws, $row, $time, $totalLenHex, $idHex, $actualLenHex, $checksumHex, $totalLenDec, $idDec, $actualLenDec, $checksumDec)

    $ws.Cells.Item($row, 1).Value = $time
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $totalLenHex
    $ws.Cells.Item($row, 3).Value = $idHex
    $ws.Cells.Item($row, 4).Value = $actualLenHex
    $ws.Cells.Item($row, 5).Value = $checksumHex
    $ws.Cells.Item($row, 6).Value = $totalLenDec
    $ws.Cells.Item($row, 7).Value = $idDec
    $ws.Cells.Item($row, 8).Value = $actualLenDec
    $ws.Cells.Item($row, 9).Value = $checksumDec
}

$wb = $excel.ActiveWorkbook

# ROW35-FE-LIFTER
$ws1 = $wb.Worksheets.Item("ROW35-FE-LIFTER")
$idDec1 = [double]"5.68631262647114e+23"
Add-UploadRow $ws1 46 45747.33911065972 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x76" "0xd" 400 $idDec1 374 13

# ROW35-MID-LIFTER
$ws2 = $wb.Worksheets.Item("ROW35-MID-LIFTER")
$idDec2 = [double]"5.68631262647114e+23"
Add-UploadRow $ws2 46 45747.18948653936 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x76" "0xe" 400 $idDec2 374 14

# ROW02-FE-LIFTER
$ws3 = $wb.Worksheets.Item("ROW02-FE-LIFTER")
$idDec3 = [double]"5.68631262647114e+23"
Add-UploadRow $ws3 46 45747.33085461806 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x76" "0x3" 400 $idDec3 374 3

# ROW02-MID-LIFTER
$ws4 = $wb.Worksheets.Item("ROW02-MID-LIFTER")
$idDec4 = [double]"9.85046333984776e+23"
Add-UploadRow $ws4 46 45747.38666480324 "0x01,0x90" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x01,0x76" "0x3" 400 $idDec4 374 3
